$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: update progress percentages in column I ---
$ws.Range("I7").Value = 0.7
$ws.Range("I8").Value = 0.7
$ws.Range("I10").Value = 0.7
$ws.Range("I12").Value = 0.7
$ws.Range("I16").Value = 0.7
$ws.Range("I18").Value = 0.6
$ws.Range("I19").Value = 0.7
$ws.Range("I21").Value = 0.5
$ws.Range("I24").Value = 0.6
$ws.Range("I25").Value = 0.6
$ws.Range("I27").Value = 0.4
$ws.Range("I29").Value = 0.4
$ws.Range("I30").Value = 0.4
$ws.Range("I31").Value = 0.4
$ws.Range("I36").Value = 0.7
$ws.Range("I37").Value = 0.7
$ws.Range("I38").Value = 0.4
$ws.Range("I47").Value = 1
$ws.Range("I48").Value = 0.5
$ws.Range("I51").Value = 0.7
$ws.Range("I52").Value = 0.7

# --- Step 2: stage the 7 distinct B-column formats used as swap sources ---
# (capture them in scratch cells before any destination is overwritten)
$ws.Range("B3").Copy() | Out-Null
$ws.Range("ZZ100").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("ZZ101").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Copy() | Out-Null
$ws.Range("ZZ102").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("ZZ103").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Copy() | Out-Null
$ws.Range("ZZ104").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Copy() | Out-Null
$ws.Range("ZZ105").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Copy() | Out-Null
$ws.Range("ZZ106").PasteSpecial(-4122) | Out-Null

# --- Step 3: apply the swapped formats onto the B column cells (rows 3-52) ---
$ws.Range("ZZ103").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ104").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ105").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ106").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ100").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ100").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ100").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ100").Copy() | Out-Null
$ws.Range("B34").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B38").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ100").Copy() | Out-Null
$ws.Range("B39").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B40").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B41").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B42").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B43").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B44").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B45").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B46").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B47").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B48").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B49").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ100").Copy() | Out-Null
$ws.Range("B50").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ101").Copy() | Out-Null
$ws.Range("B51").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ102").Copy() | Out-Null
$ws.Range("B52").PasteSpecial(-4122) | Out-Null

# --- Step 4: clean up scratch cells so the sheet dimension is unaffected ---
$ws.Range("ZZ100").Clear() | Out-Null
$ws.Range("ZZ101").Clear() | Out-Null
$ws.Range("ZZ102").Clear() | Out-Null
$ws.Range("ZZ103").Clear() | Out-Null
$ws.Range("ZZ104").Clear() | Out-Null
$ws.Range("ZZ105").Clear() | Out-Null
$ws.Range("ZZ106").Clear() | Out-Null

$excel.CutCopyMode = 0

# --- Step 5: update sheet view (scroll position + selection) ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("N12").Select() | Out-Null
